$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "JSU(-0.8415878686863555, 1.0101615247787215, 0.6585763923496767, 2.499121682956182)"
$ws.Range("C2").Value = "NIG(0.8258016321049033, 0.5857586405327566, 4.497054452154432, 4.8708672604268655)"
$ws.Range("D2").Value = "MIE(6.877818745104607, 4.297349092046274, -7.093176424555315, 9.215903763746997)"
$ws.Range("E2").Value = "JSU(-1.4838472377302376, 1.2600346295724658, 1.5466732329853992, 5.524938060679077)"
